# Adds two new columns, I ("I0") and J ("IF"), to the active worksheet.
# Header row (row 1) gets the labels "I0" and "IF"; rows 2-74 get the
# corresponding numeric values that were added alongside the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - match the formatting already used by the other header cells
# (bold font, thin border, centered alignment) by copying it from H1, then
# set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I74 and J2:J74
$iValues = @(7,5,7,7,9,6,6,8,7,7,5,6,10,7,5,7,6,7,7,9,8,5,8,7,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,8,9,8,9,9,8,8,9,9,9,9,9,9,8,9,8,8,9,7,8,7,7,7,7,8,8,7,8,6,9,3,2,4)
$jValues = @(7,5,7,7,9,6,6,8,7,7,5,6,10,7,6,7,6,8,7,9,8,6,8,7,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,8,9,8,9,9,9,8,9,9,9,9,9,9,8,9,8,8,9,7,8,7,7,7,7,8,8,7,8,7,9,3,2,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
